$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Services")
$cell = $ws.Cells.Item(57, 13)
$ws.Hyperlinks.Add($cell, "https://www.ccac.edu/workforce/index.php")
"done"
